$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117; this pushes the existing rows 117:190 down to 118:191
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new weekly data entry
$ws.Range("A117").Value = 5
$ws.Range("B117").Value = "Macroferia Regional de Talca"
$ws.Range("C117").Value = "Maule"
$ws.Range("D117").Value = 44907
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 100112031
$ws.Range("G117").Value = "Poroto verde"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("M117").Value = 30000
$ws.Range("N117").Value = "$/saco 25 kilos"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 1200
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"

Write-Host "Row 117 inserted and populated"
